$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet: insert a new (blank) column before
# column N, shifting the existing "Late", "heading" and "Outstanding" columns
# one place to the right.
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Select the sheet and a specific cell, then activate the sheet so it
# becomes the workbook's active tab.
$ws.Select() | Out-Null
$ws.Range("P5").Select() | Out-Null
$ws.Activate() | Out-Null
